$d = $word.ActiveDocument

# The edit simplifies <w:docDefaults> in word/styles.xml: drop the
# explicit "value == schema default" properties that Word's exporter had
# spelled out, keeping only the ones that actually diverge from the
# built-in defaults. We do this by round-tripping the package through
# Document.WordOpenXML (which embeds the real part XML verbatim) so the
# change lands precisely on styles.xml's docDefaults element.

$xml = $d.WordOpenXML

$old = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:val="clear" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:pBdr><w:top w:val="nil" w:sz="0" w:space="0"/><w:left w:val="nil" w:sz="0" w:space="0"/><w:bottom w:val="nil" w:sz="0" w:space="0"/><w:right w:val="nil" w:sz="0" w:space="0"/><w:between w:val="nil" w:sz="0" w:space="0"/></w:pBdr><w:shd w:val="clear" w:fill="auto"/><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr></w:pPrDefault></w:docDefaults>'

$new = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

if ($xml.Contains($old)) {
    $xml = $xml.Replace($old, $new)
    $d.WordOpenXML = $xml
    Write-Output "docDefaults simplified"
} else {
    Write-Output "WARNING: expected docDefaults block not found"
}
